# Auto-generated script applying targeted cell value updates
# as described by the commit diff (numeric recalculations in
# the profit-tracking tables of each job sheet).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2566.6667
$ws.Range("I112").Value = 333.33334
$ws.Range("J112").Value = 2790
$ws.Range("K112").Value = 1000.00002
$ws.Range("L112").Value = 8370
$ws.Range("M112").Value = 107.9999799999999
$ws.Range("N112").Value = -10586
$ws.Range("H137").Value = 2632614.5
$ws.Range("I137").Value = 1137364.4
$ws.Range("J137").Value = 7693461.5
$ws.Range("K137").Value = 3412093.2
$ws.Range("L137").Value = 23080384.5
$ws.Range("M137").Value = -3409543.2
$ws.Range("N137").Value = -23085484.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1391.2307
$ws.Range("I61").Value = 1454.8667
$ws.Range("J61").Value = 1179.1111
$ws.Range("K61").Value = 1454.8667
$ws.Range("L61").Value = 1179.1111
$ws.Range("M61").Value = -1242.8667
$ws.Range("N61").Value = -1603.1111
$ws.Range("H74").Value = 883.33844
$ws.Range("I74").Value = 855.9259
$ws.Range("J74").Value = 1017.9091
$ws.Range("K74").Value = 855.9259
$ws.Range("L74").Value = 1017.9091
$ws.Range("M74").Value = 18.07410000000004
$ws.Range("N74").Value = -2765.9091
$ws.Range("H77").Value = 883.33844
$ws.Range("I77").Value = 855.9259
$ws.Range("J77").Value = 1017.9091
$ws.Range("K77").Value = 4279.6295
$ws.Range("L77").Value = 5089.5455
$ws.Range("M77").Value = 88.37049999999999
$ws.Range("N77").Value = -13825.5455
$ws.Range("H102").Value = 166668500
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 333335000
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 333335000
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -333338244
$ws.Range("H136").Value = 1391.2307
$ws.Range("I136").Value = 1454.8667
$ws.Range("J136").Value = 1179.1111
$ws.Range("K136").Value = 4364.6001
$ws.Range("L136").Value = 3537.3333
$ws.Range("M136").Value = -1814.6001
$ws.Range("N136").Value = -8637.3333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 104518.34
$ws.Range("I134").Value = 131382.27
$ws.Range("J134").Value = 1540
$ws.Range("K134").Value = 394146.8099999999
$ws.Range("L134").Value = 4620
$ws.Range("M134").Value = -391611.8099999999
$ws.Range("N134").Value = -9690

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2250.2058
$ws.Range("I31").Value = 1663.963
$ws.Range("J31").Value = 4511.4287
$ws.Range("K31").Value = 1663.963
$ws.Range("L31").Value = 4511.4287
$ws.Range("M31").Value = -1368.963
$ws.Range("N31").Value = -5101.4287
$ws.Range("H34").Value = 2250.2058
$ws.Range("I34").Value = 1663.963
$ws.Range("J34").Value = 4511.4287
$ws.Range("K34").Value = 1663.963
$ws.Range("L34").Value = 4511.4287
$ws.Range("M34").Value = -1461.963
$ws.Range("N34").Value = -4915.4287
$ws.Range("H58").Value = 1589.3572
$ws.Range("I58").Value = 1765.9714
$ws.Range("J58").Value = 706.2857
$ws.Range("K58").Value = 1765.9714
$ws.Range("L58").Value = 706.2857
$ws.Range("M58").Value = -1562.9714
$ws.Range("N58").Value = -1112.2857
$ws.Range("H132").Value = 2901.325
$ws.Range("I132").Value = 2478.6775
$ws.Range("J132").Value = 4357.1113
$ws.Range("K132").Value = 7436.032499999999
$ws.Range("L132").Value = 13071.3339
$ws.Range("M132").Value = -4906.032499999999
$ws.Range("N132").Value = -18131.3339
$ws.Range("H134").Value = 5327.4883
$ws.Range("I134").Value = 5868.3887
$ws.Range("K134").Value = 17605.1661
$ws.Range("M134").Value = -15070.1661
$ws.Range("H136").Value = 1589.3572
$ws.Range("I136").Value = 1765.9714
$ws.Range("J136").Value = 706.2857
$ws.Range("K136").Value = 5297.914199999999
$ws.Range("L136").Value = 2118.8571
$ws.Range("M136").Value = -2747.914199999999
$ws.Range("N136").Value = -7218.8571

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 193508.78
$ws.Range("I5").Value = 222.7541
$ws.Range("J5").Value = 1667314.8
$ws.Range("K5").Value = 668.2623
$ws.Range("L5").Value = 5001944.4
$ws.Range("M5").Value = -556.2623
$ws.Range("N5").Value = -5002168.4
$ws.Range("H122").Value = 50895.703
$ws.Range("I122").Value = 304.15384
$ws.Range("J122").Value = 60158.945
$ws.Range("K122").Value = 2737.38456
$ws.Range("L122").Value = 541430.505
$ws.Range("M122").Value = -287.38456
$ws.Range("N122").Value = -546330.505
$ws.Range("H129").Value = 1448.421
$ws.Range("I129").Value = 760.7692
$ws.Range("J129").Value = 2938.3333
$ws.Range("K129").Value = 2282.3076
$ws.Range("L129").Value = 8814.999899999999
$ws.Range("M129").Value = 2717.6924
$ws.Range("N129").Value = -18814.9999
$ws.Range("H130").Value = 1889.091
$ws.Range("I130").Value = 780
$ws.Range("K130").Value = 2340
$ws.Range("M130").Value = 2680
$ws.Range("H131").Value = 933.75
$ws.Range("I131").Value = 537.1429000000001
$ws.Range("J131").Value = 963.6022
$ws.Range("K131").Value = 1611.4287
$ws.Range("L131").Value = 2890.8066
$ws.Range("M131").Value = 3428.5713
$ws.Range("N131").Value = -12970.8066
$ws.Range("H135").Value = 193508.78
$ws.Range("I135").Value = 222.7541
$ws.Range("J135").Value = 1667314.8
$ws.Range("K135").Value = 2004.7869
$ws.Range("L135").Value = 15005833.2
$ws.Range("M135").Value = 530.2130999999999
$ws.Range("N135").Value = -15010903.2
$ws.Range("H136").Value = 3728
$ws.Range("I136").Value = 1740.5555
$ws.Range("J136").Value = 4390.4814
$ws.Range("K136").Value = 5221.666499999999
$ws.Range("L136").Value = 13171.4442
$ws.Range("M136").Value = -121.6664999999994
$ws.Range("N136").Value = -23371.4442
$ws.Range("H139").Value = 23316.234
$ws.Range("J139").Value = 103390.3
$ws.Range("L139").Value = 310170.9
$ws.Range("N139").Value = -320450.9
$ws.Range("H140").Value = 42567.617
$ws.Range("I140").Value = 144441.42
$ws.Range("J140").Value = 5035.1577
$ws.Range("K140").Value = 433324.26
$ws.Range("L140").Value = 15105.4731
$ws.Range("M140").Value = -428144.26
$ws.Range("N140").Value = -25465.4731

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1671.6271
$ws.Range("I132").Value = 1600.551
$ws.Range("K132").Value = 4801.653
$ws.Range("M132").Value = -2271.653
$ws.Range("H136").Value = 1385.2373
$ws.Range("I136").Value = 1288.3265
$ws.Range("J136").Value = 1860.1
$ws.Range("K136").Value = 3864.979499999999
$ws.Range("L136").Value = 5580.299999999999
$ws.Range("M136").Value = -1314.979499999999
$ws.Range("N136").Value = -10680.3

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1698.0526
$ws.Range("I132").Value = 2134.8462
$ws.Range("K132").Value = 6404.5386
$ws.Range("M132").Value = -3874.5386
$ws.Range("H136").Value = 1582.9445
$ws.Range("I136").Value = 1692.0714
$ws.Range("J136").Value = 1201
$ws.Range("K136").Value = 5076.2142
$ws.Range("L136").Value = 3603
$ws.Range("M136").Value = -2526.2142
$ws.Range("N136").Value = -8703
